$d = $word.ActiveDocument

# The last paragraph in the document is currently empty (it only holds the
# _GoBack bookmark). Place the insertion point at the very end of the
# document, type "Four" there, then press Enter to start a new (empty)
# paragraph after it - mirroring how a user would continue the "One/Two/
# Three/..." list.
$end = $d.Content
$end.Collapse(0)                 # wdCollapseEnd
$end.InsertAfter("Four")
$end.Collapse(0)                 # move past the inserted text
$end.InsertParagraphAfter()
